$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the Kaizer Chiefs vs Orlando Pirates / Polokwane vs Swallows rows ---
# Row 83 currently holds "Kaizer Chiefs vs Orlando Pirates" data and row 84
# holds "Polokwane vs Swallows" data. The correct order is Polokwane/Swallows
# first (row 83) then Kaizer Chiefs/Orlando Pirates (row 84). Columns A:E
# (index/country/tournament/season/date) stay put; only F:V swap.

$row83 = @("Kaizer Chiefs", 0, "Orlando Pirates", 1, 2.79, "11/11/2023 07:48", 3.17, "11/11/2023 14:21", 2.84, "11/11/2023 07:48", 2.87, "11/11/2023 14:21", 2.79, "11/11/2023 07:48", 2.64, "11/11/2023 14:21", "https://www.betexplorer.com/football/south-africa/premier-league/kaizer-chiefs-orlando-pirates/0MutWbLr/")
$row84 = @("Polokwane", 0, "Swallows", 0, 2.65, "11/11/2023 07:48", 2.76, "11/11/2023 14:21", 2.86, "11/11/2023 07:48", 2.84, "11/11/2023 14:21", 2.92, "11/11/2023 07:48", 3.04, "11/11/2023 14:21", "https://www.betexplorer.com/football/south-africa/premier-league/polokwane-city-swallows-fc/fqupVIzk/")

# Write row84's data (Polokwane/Swallows) into row 83
for ($i = 0; $i -lt $row84.Length; $i++) {
    $ws.Cells.Item(83, 6 + $i).Value = $row84[$i]
}

# Write row83's data (Kaizer Chiefs/Orlando Pirates) into row 84
for ($i = 0; $i -lt $row83.Length; $i++) {
    $ws.Cells.Item(84, 6 + $i).Value = $row83[$i]
}

# --- Append the new AmaZulu vs Cape Town Spurs match as row 109 ---
# Copy the formatting (styles) of row 108's index/date cells so the new
# row matches the sheet's existing look (bold/bordered index column,
# date-formatted match-date column).
$ws.Cells.Item(108, 1).Copy($ws.Cells.Item(109, 1))
$ws.Cells.Item(108, 5).Copy($ws.Cells.Item(109, 5))

$ws.Cells.Item(109, 1).Value = 108
$ws.Cells.Item(109, 2).Value = "south-africa"
$ws.Cells.Item(109, 3).Value = "premier-league"
$ws.Cells.Item(109, 4).Value = "2023-2024"
$ws.Cells.Item(109, 5).Value = 45282.77083333334
$ws.Cells.Item(109, 6).Value = "AmaZulu"
$ws.Cells.Item(109, 7).Value = 1
$ws.Cells.Item(109, 8).Value = "Cape Town Spurs"
$ws.Cells.Item(109, 9).Value = 1
$ws.Cells.Item(109, 10).Value = 1.64
$ws.Cells.Item(109, 11).Value = "20/12/2023 03:12"
$ws.Cells.Item(109, 12).Value = 1.75
$ws.Cells.Item(109, 13).Value = "22/12/2023 17:54"
$ws.Cells.Item(109, 14).Value = 3.27
$ws.Cells.Item(109, 15).Value = "20/12/2023 03:12"
$ws.Cells.Item(109, 16).Value = 3.25
$ws.Cells.Item(109, 17).Value = "22/12/2023 18:01"
$ws.Cells.Item(109, 18).Value = 6.03
$ws.Cells.Item(109, 19).Value = "20/12/2023 03:12"
$ws.Cells.Item(109, 20).Value = 5.88
$ws.Cells.Item(109, 21).Value = "22/12/2023 18:01"
$ws.Cells.Item(109, 22).Value = "https://www.betexplorer.com/football/south-africa/premier-league/amazulu-cape-town-spurs/zZjEiXmJ/"
